$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.200255632400513
$ws.Range("B1").Value = 1.168312907218933
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.033358097076416
$ws.Range("E1").Value = 0.9699259996414185
